# Auto-generated Excel COM-interop edit script
# Applies numeric updates to H:N columns (price/profit data) across several
# sheets of the workbook, matching the target OOXML diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 4200
$ws.Range("I21").Value = 4200
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 4200
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -3732
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 4200
$ws.Range("I23").Value = 4200
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 4200
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -3966
$ws.Range("N23").ClearContents()

$ws.Range("H37").Value = 1352.6666
$ws.Range("I37").Value = 529
$ws.Range("J37").Value = 3000
$ws.Range("K37").Value = 1587
$ws.Range("L37").Value = 9000
$ws.Range("M37").Value = -1461
$ws.Range("N37").Value = -9252

$ws.Range("H70").Value = 7257.6665
$ws.Range("I70").Value = 4122.25
$ws.Range("J70").Value = 9766
$ws.Range("K70").Value = 12366.75
$ws.Range("L70").Value = 29298
$ws.Range("M70").Value = -12096.75
$ws.Range("N70").Value = -29838

$ws.Range("H73").Value = 7257.6665
$ws.Range("I73").Value = 4122.25
$ws.Range("J73").Value = 9766
$ws.Range("K73").Value = 12366.75
$ws.Range("L73").Value = 29298
$ws.Range("M73").Value = -11430.75
$ws.Range("N73").Value = -31170

$ws.Range("H141").Value = 997.3125
$ws.Range("I141").Value = 905.5714
$ws.Range("J141").Value = 1639.5
$ws.Range("K141").Value = 2716.7142
$ws.Range("L141").Value = 4918.5
$ws.Range("M141").Value = 2463.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 1270.1666
$ws.Range("I50").Value = 574
$ws.Range("J50").Value = 1618.25
$ws.Range("K50").Value = 574
$ws.Range("L50").Value = 1618.25
$ws.Range("M50").Value = 140
$ws.Range("N50").Value = -3046.25

$ws.Range("H132").Value = 3488.9666
$ws.Range("I132").Value = 3109.1482
$ws.Range("J132").Value = 6907.3335
$ws.Range("K132").Value = 9327.444600000001
$ws.Range("L132").Value = 20722.0005
$ws.Range("M132").Value = -6797.444600000001
$ws.Range("N132").Value = -25782.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 52095.75
$ws.Range("I80").Value = 357.83334
$ws.Range("J80").Value = 74269.14
$ws.Range("K80").Value = 357.83334
$ws.Range("L80").Value = 74269.14
$ws.Range("M80").Value = 640.16666
$ws.Range("N80").Value = -76265.14

$ws.Range("H83").Value = 52095.75
$ws.Range("I83").Value = 357.83334
$ws.Range("J83").Value = 74269.14
$ws.Range("K83").Value = 1789.1667
$ws.Range("L83").Value = 371345.7
$ws.Range("M83").Value = 3202.8333
$ws.Range("N83").Value = -381329.7

$ws.Range("H107").Value = 1341.6666
$ws.Range("I107").Value = 1260
$ws.Range("J107").Value = 1535.625
$ws.Range("K107").Value = 1260
$ws.Range("L107").Value = 1535.625
$ws.Range("M107").Value = 660
$ws.Range("N107").Value = -5375.625

$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H134").Value = 4576.6206
$ws.Range("I134").Value = 864.8889
$ws.Range("J134").Value = 10650.363
$ws.Range("K134").Value = 2594.6667
$ws.Range("L134").Value = 31951.089
$ws.Range("M134").Value = -59.66670000000022

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 2321.5
$ws.Range("I8").Value = 643
$ws.Range("J8").Value = 4000
$ws.Range("K8").Value = 643
$ws.Range("L8").Value = 4000
$ws.Range("M8").Value = -503

$ws.Range("H70").Value = 79999
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 79999
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 79999
$ws.Range("N70").Value = -80629

$ws.Range("H73").Value = 79999
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 79999
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 79999
$ws.Range("N73").Value = -82183

$ws.Range("H94").Value = 5965.2
$ws.Range("I94").Value = 9116.333000000001
$ws.Range("J94").Value = 1238.5
$ws.Range("K94").Value = 9116.333000000001
$ws.Range("L94").Value = 1238.5
$ws.Range("M94").Value = -8665.333000000001
$ws.Range("N94").Value = -2140.5

$ws.Range("H124").Value = 1218333
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 1218333
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 1218333
$ws.Range("N124").Value = -1223243

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3567218.2
$ws.Range("I4").Value = 2746340.2
$ws.Range("J4").Value = 5893039
$ws.Range("K4").Value = 8239020.600000001
$ws.Range("L4").Value = 17679117
$ws.Range("M4").Value = -8238908.600000001

$ws.Range("H14").Value = 1082.1538
$ws.Range("I14").Value = 1082.1538
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 3246.4614
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -3073.4614

$ws.Range("H74").Value = 7050
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 7050
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 21150
$ws.Range("N74").Value = -23272
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 7050
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 7050
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 63450
$ws.Range("N77").Value = -74058
$ws.Range("M77").ClearContents()

$ws.Range("H131").Value = 13895726
$ws.Range("I131").Value = 33333940
$ws.Range("J131").Value = 11287.857
$ws.Range("K131").Value = 100001820
$ws.Range("L131").Value = 33863.571
$ws.Range("M131").Value = -99996780
$ws.Range("N131").Value = -43943.571

$ws.Range("H140").Value = 1250
$ws.Range("I140").Value = 1000
$ws.Range("J140").Value = 2000
$ws.Range("K140").Value = 3000
$ws.Range("L140").Value = 6000
$ws.Range("M140").Value = 2180
$ws.Range("N140").Value = -16360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6968.933
$ws.Range("I70").Value = 3870.375
$ws.Range("J70").Value = 10510.143
$ws.Range("K70").Value = 3870.375
$ws.Range("L70").Value = 10510.143
$ws.Range("M70").Value = -3600.375
$ws.Range("N70").Value = -11050.143

$ws.Range("H73").Value = 6968.933
$ws.Range("I73").Value = 3870.375
$ws.Range("J73").Value = 10510.143
$ws.Range("K73").Value = 3870.375
$ws.Range("L73").Value = 10510.143
$ws.Range("M73").Value = -2934.375
$ws.Range("N73").Value = -12382.143

$ws.Range("H132").Value = 7515.905
$ws.Range("I132").Value = 6018.6665
$ws.Range("J132").Value = 16499.334
$ws.Range("K132").Value = 18055.9995
$ws.Range("L132").Value = 49498.00199999999
$ws.Range("M132").Value = -15525.9995
$ws.Range("N132").Value = -54558.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H21").Value = 5000
$ws.Range("I21").Value = 5000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -4826
$ws.Range("N21").ClearContents()

$ws.Range("H93").Value = 1859.1428
$ws.Range("I93").Value = 1802
$ws.Range("J93").Value = 2002
$ws.Range("K93").Value = 1802
$ws.Range("L93").Value = 2002
$ws.Range("M93").Value = -554

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1319.2
$ws.Range("I14").Value = 1319.2
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1319.2
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1151.2
$ws.Range("N14").ClearContents()

$ws.Range("H20").Value = 45597.75
$ws.Range("I20").Value = 12500
$ws.Range("J20").Value = 78695.5
$ws.Range("K20").Value = 12500
$ws.Range("L20").Value = 78695.5
$ws.Range("M20").Value = -12260
$ws.Range("N20").Value = -79175.5

$ws.Range("H38").Value = 31000
$ws.Range("I38").Value = 40000
$ws.Range("J38").Value = 22000
$ws.Range("K38").Value = 40000
$ws.Range("L38").Value = 22000
$ws.Range("M38").Value = -39527
$ws.Range("N38").Value = -22946
